$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws = $wb.Worksheets.Item("Summary")

# Insert a new column at N (14) to make room for "Files Count";
# this shifts the existing Analysis Date (N) -> O and Data Date (O) -> P.
$ws.Columns.Item(14).Insert()

# New header for the inserted column
$ws.Range("N1").Value = "Files Count"

# Force the Annual Return column to stay text (avoid Excel auto-% conversion)
$ws.Range("G2:G16").NumberFormat = "@"

$ws.Range("G2").Value = "+216.99%"
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 6
$ws.Range("G3").Value = "+195.90%"
$ws.Range("M3").Value = 6
$ws.Range("N3").Value = 6
$ws.Range("G4").Value = "+51.47%"
$ws.Range("M4").Value = 6
$ws.Range("N4").Value = 6
$ws.Range("G5").Value = "+33.17%"
$ws.Range("M5").Value = 6
$ws.Range("N5").Value = 6
$ws.Range("G6").Value = "+3.27%"
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 6
$ws.Range("G7").Value = "+17.36%"
$ws.Range("M7").Value = 6
$ws.Range("N7").Value = 6
$ws.Range("G8").Value = "+108.09%"
$ws.Range("M8").Value = 6
$ws.Range("N8").Value = 6
$ws.Range("G9").Value = "+17.68%"
$ws.Range("M9").Value = 6
$ws.Range("N9").Value = 6
$ws.Range("G10").Value = "+5.95%"
$ws.Range("M10").Value = 6
$ws.Range("N10").Value = 6
$ws.Range("G11").Value = "+2.14%"
$ws.Range("M11").Value = 6
$ws.Range("N11").Value = 6
$ws.Range("G12").Value = "+16.98%"
$ws.Range("M12").Value = 5
$ws.Range("N12").Value = 5
$ws.Range("G13").Value = "+134.88%"
$ws.Range("M13").Value = 6
$ws.Range("N13").Value = 6
$ws.Range("G14").Value = "+126.49%"
$ws.Range("M14").Value = 6
$ws.Range("N14").Value = 6
$ws.Range("G15").Value = "+3.66%"
$ws.Range("M15").Value = 6
$ws.Range("N15").Value = 6
$ws.Range("G16").Value = "+39.58%"
$ws.Range("M16").Value = 6
$ws.Range("N16").Value = 6

# ---- Sheet: Pattern1-Pure Data ----
$ws = $wb.Worksheets.Item("Pattern1-Pure Data")

# Insert a new column at N (14) to make room for "Files Count";
# this shifts the existing Analysis Date (N) -> O and Data Date (O) -> P.
$ws.Columns.Item(14).Insert()

# New header for the inserted column
$ws.Range("N1").Value = "Files Count"

# Force the Annual Return column to stay text (avoid Excel auto-% conversion)
$ws.Range("G2:G6").NumberFormat = "@"

$ws.Range("G2").Value = "+216.99%"
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 6
$ws.Range("G3").Value = "+195.90%"
$ws.Range("M3").Value = 6
$ws.Range("N3").Value = 6
$ws.Range("G4").Value = "+51.47%"
$ws.Range("M4").Value = 6
$ws.Range("N4").Value = 6
$ws.Range("G5").Value = "+33.17%"
$ws.Range("M5").Value = 6
$ws.Range("N5").Value = 6
$ws.Range("G6").Value = "+3.27%"
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 6

# ---- Sheet: Pattern2-Data+Technical ----
$ws = $wb.Worksheets.Item("Pattern2-Data+Technical")

# Insert a new column at N (14) to make room for "Files Count";
# this shifts the existing Analysis Date (N) -> O and Data Date (O) -> P.
$ws.Columns.Item(14).Insert()

# New header for the inserted column
$ws.Range("N1").Value = "Files Count"

# Force the Annual Return column to stay text (avoid Excel auto-% conversion)
$ws.Range("G2:G6").NumberFormat = "@"

$ws.Range("G2").Value = "+17.36%"
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 6
$ws.Range("G3").Value = "+108.09%"
$ws.Range("M3").Value = 6
$ws.Range("N3").Value = 6
$ws.Range("G4").Value = "+17.68%"
$ws.Range("M4").Value = 6
$ws.Range("N4").Value = 6
$ws.Range("G5").Value = "+5.95%"
$ws.Range("M5").Value = 6
$ws.Range("N5").Value = 6
$ws.Range("G6").Value = "+2.14%"
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 6

# ---- Sheet: Pattern3-Data+News ----
$ws = $wb.Worksheets.Item("Pattern3-Data+News")

# Insert a new column at N (14) to make room for "Files Count";
# this shifts the existing Analysis Date (N) -> O and Data Date (O) -> P.
$ws.Columns.Item(14).Insert()

# New header for the inserted column
$ws.Range("N1").Value = "Files Count"

# Force the Annual Return column to stay text (avoid Excel auto-% conversion)
$ws.Range("G2:G6").NumberFormat = "@"

$ws.Range("G2").Value = "+16.98%"
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 5
$ws.Range("G3").Value = "+134.88%"
$ws.Range("M3").Value = 6
$ws.Range("N3").Value = 6
$ws.Range("G4").Value = "+126.49%"
$ws.Range("M4").Value = 6
$ws.Range("N4").Value = 6
$ws.Range("G5").Value = "+3.66%"
$ws.Range("M5").Value = 6
$ws.Range("N5").Value = 6
$ws.Range("G6").Value = "+39.58%"
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 6
